$wb = $excel.ActiveWorkbook

# Sheet ALC, Row 32
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(32, 8).Value = 420.01514  # H32: 515.29034 -> 420.01514
$ws.Cells.Item(32, 9).Value = 513.2727  # I32: 522.36365 -> 513.2727
$ws.Cells.Item(32, 10).Value = 326.75757  # J32: 507.2414 -> 326.75757
$ws.Cells.Item(32, 11).Value = 513.2727  # K32: 522.36365 -> 513.2727
$ws.Cells.Item(32, 12).Value = 326.75757  # L32: 507.2414 -> 326.75757
$ws.Cells.Item(32, 13).Value = -187.2727  # M32: -196.36365 -> -187.2727
$ws.Cells.Item(32, 14).Value = -978.75757  # N32: -1159.2414 -> -978.75757

# Sheet ALC, Row 41
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(41, 8).Value = 606.36365  # H41: 150.5 -> 606.36365
$ws.Cells.Item(41, 9).Value = 557  # I41: 150.5 -> 557
$ws.Cells.Item(41, 10).Value = 647.5  # J41: 0 -> 647.5
$ws.Cells.Item(41, 11).Value = 557  # K41: 150.5 -> 557
$ws.Cells.Item(41, 12).Value = 647.5  # L41: 0 -> 647.5
$ws.Cells.Item(41, 13).Value = -117  # M41: 289.5 -> -117
$ws.Cells.Item(41, 14).Value = -1527.5  # N41: None -> -1527.5

# Sheet ALC, Row 98
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(98, 8).Value = 35298.934  # H98: 79377.30499999999 -> 35298.934
$ws.Cells.Item(98, 9).Value = 38932.15  # I98: 85742.086 -> 38932.15
$ws.Cells.Item(98, 10).Value = 2600  # J98: 3000 -> 2600
$ws.Cells.Item(98, 11).Value = 38932.15  # K98: 85742.086 -> 38932.15
$ws.Cells.Item(98, 12).Value = 2600  # L98: 3000 -> 2600
$ws.Cells.Item(98, 13).Value = -37434.15  # M98: -84244.086 -> -37434.15
$ws.Cells.Item(98, 14).Value = -5596  # N98: -5996 -> -5596

# Sheet ALC, Row 122
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(122, 8).Value = 35298.934  # H122: 79377.30499999999 -> 35298.934
$ws.Cells.Item(122, 9).Value = 38932.15  # I122: 85742.086 -> 38932.15
$ws.Cells.Item(122, 10).Value = 2600  # J122: 3000 -> 2600
$ws.Cells.Item(122, 11).Value = 116796.45  # K122: 257226.258 -> 116796.45
$ws.Cells.Item(122, 12).Value = 7800  # L122: 9000 -> 7800
$ws.Cells.Item(122, 13).Value = -114346.45  # M122: -254776.258 -> -114346.45
$ws.Cells.Item(122, 14).Value = -12700  # N122: -13900 -> -12700

# Sheet ALC, Row 131
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(131, 8).Value = 1655.7142  # H131: 1146.875 -> 1655.7142
$ws.Cells.Item(131, 9).Value = 1018  # I131: 810.7143 -> 1018
$ws.Cells.Item(131, 10).Value = 3250  # J131: 3500 -> 3250
$ws.Cells.Item(131, 11).Value = 3054  # K131: 2432.1429 -> 3054
$ws.Cells.Item(131, 12).Value = 9750  # L131: 10500 -> 9750
$ws.Cells.Item(131, 13).Value = 1986  # M131: 2607.8571 -> 1986
$ws.Cells.Item(131, 14).Value = -19830  # N131: -20580 -> -19830

# Sheet ALC, Row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(132, 8).Value = 3970397.5  # H132: 1066940.4 -> 3970397.5
$ws.Cells.Item(132, 9).Value = 4465947  # I132: 1099600.1 -> 4465947
$ws.Cells.Item(132, 10).Value = 6000  # J132: 5500 -> 6000
$ws.Cells.Item(132, 11).Value = 13397841  # K132: 3298800.3 -> 13397841
$ws.Cells.Item(132, 12).Value = 18000  # L132: 16500 -> 18000
$ws.Cells.Item(132, 13).Value = -13395311  # M132: -3296270.3 -> -13395311
$ws.Cells.Item(132, 14).Value = -23060  # N132: -21560 -> -23060

# Sheet ALC, Row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(137, 8).Value = 1159.1025  # H137: 1171.6154 -> 1159.1025
$ws.Cells.Item(137, 9).Value = 953.90625  # I137: 967.1 -> 953.90625
$ws.Cells.Item(137, 10).Value = 2097.1428  # J137: 1853.3334 -> 2097.1428
$ws.Cells.Item(137, 11).Value = 2861.71875  # K137: 2901.3 -> 2861.71875
$ws.Cells.Item(137, 12).Value = 6291.428400000001  # L137: 5560.0002 -> 6291.428400000001
$ws.Cells.Item(137, 13).Value = -311.71875  # M137: -351.3000000000002 -> -311.71875
$ws.Cells.Item(137, 14).Value = -11391.4284  # N137: -10660.0002 -> -11391.4284

# Sheet ARM, Row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 1528.5927  # H61: 1077.1818 -> 1528.5927
$ws.Cells.Item(61, 9).Value = 1117.95  # I61: 729.72974 -> 1117.95
$ws.Cells.Item(61, 10).Value = 2701.8572  # J61: 2913.7144 -> 2701.8572
$ws.Cells.Item(61, 11).Value = 1117.95  # K61: 729.72974 -> 1117.95
$ws.Cells.Item(61, 12).Value = 2701.8572  # L61: 2913.7144 -> 2701.8572
$ws.Cells.Item(61, 13).Value = -905.95  # M61: -517.72974 -> -905.95
$ws.Cells.Item(61, 14).Value = -3125.8572  # N61: -3337.7144 -> -3125.8572

# Sheet ARM, Row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(74, 8).Value = 805.8570999999999  # H74: 838.5789 -> 805.8570999999999
$ws.Cells.Item(74, 9).Value = 676.2  # I74: 704.0769 -> 676.2
$ws.Cells.Item(74, 11).Value = 676.2  # K74: 704.0769 -> 676.2
$ws.Cells.Item(74, 13).Value = 197.8  # M74: 169.9231 -> 197.8

# Sheet ARM, Row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(77, 8).Value = 805.8570999999999  # H77: 838.5789 -> 805.8570999999999
$ws.Cells.Item(77, 9).Value = 676.2  # I77: 704.0769 -> 676.2
$ws.Cells.Item(77, 11).Value = 3381  # K77: 3520.3845 -> 3381
$ws.Cells.Item(77, 13).Value = 987  # M77: 847.6154999999999 -> 987

# Sheet ARM, Row 88
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(88, 8).Value = 3057.1428  # H88: 2810 -> 3057.1428
$ws.Cells.Item(88, 9).Value = 1500  # I88: 1663.3334 -> 1500
$ws.Cells.Item(88, 10).Value = 3316.6667  # J88: 3383.3333 -> 3316.6667
$ws.Cells.Item(88, 11).Value = 1500  # K88: 1663.3334 -> 1500
$ws.Cells.Item(88, 12).Value = 3316.6667  # L88: 3383.3333 -> 3316.6667
$ws.Cells.Item(88, 13).Value = -1094  # M88: -1257.3334 -> -1094
$ws.Cells.Item(88, 14).Value = -4128.6667  # N88: -4195.3333 -> -4128.6667

# Sheet ARM, Row 91
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(91, 8).Value = 3057.1428  # H91: 2810 -> 3057.1428
$ws.Cells.Item(91, 9).Value = 1500  # I91: 1663.3334 -> 1500
$ws.Cells.Item(91, 10).Value = 3316.6667  # J91: 3383.3333 -> 3316.6667
$ws.Cells.Item(91, 11).Value = 1500  # K91: 1663.3334 -> 1500
$ws.Cells.Item(91, 12).Value = 3316.6667  # L91: 3383.3333 -> 3316.6667
$ws.Cells.Item(91, 13).Value = -96  # M91: -259.3334 -> -96
$ws.Cells.Item(91, 14).Value = -6124.6667  # N91: -6191.3333 -> -6124.6667

# Sheet ARM, Row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(122, 8).Value = 6000  # H122: 4666.6665 -> 6000
$ws.Cells.Item(122, 9).Value = 0  # I122: 5000 -> 0
$ws.Cells.Item(122, 10).Value = 6000  # J122: 4500 -> 6000
$ws.Cells.Item(122, 11).Value = 0  # K122: 15000 -> 0
$ws.Cells.Item(122, 12).Value = 18000  # L122: 13500 -> 18000
$ws.Cells.Item(122, 13).ClearContents()  # M122: was -12550
$ws.Cells.Item(122, 14).Value = -22900  # N122: -18400 -> -22900

# Sheet ARM, Row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(132, 8).Value = 1667.9667  # H132: 1433.0857 -> 1667.9667
$ws.Cells.Item(132, 9).Value = 1471.5834  # I132: 1218.2667 -> 1471.5834
$ws.Cells.Item(132, 10).Value = 2453.5  # J132: 2722 -> 2453.5
$ws.Cells.Item(132, 11).Value = 4414.7502  # K132: 3654.800099999999 -> 4414.7502
$ws.Cells.Item(132, 12).Value = 7360.5  # L132: 8166 -> 7360.5
$ws.Cells.Item(132, 13).Value = -1884.7502  # M132: -1124.800099999999 -> -1884.7502
$ws.Cells.Item(132, 14).Value = -12420.5  # N132: -13226 -> -12420.5

# Sheet ARM, Row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(136, 8).Value = 1528.5927  # H136: 1077.1818 -> 1528.5927
$ws.Cells.Item(136, 9).Value = 1117.95  # I136: 729.72974 -> 1117.95
$ws.Cells.Item(136, 10).Value = 2701.8572  # J136: 2913.7144 -> 2701.8572
$ws.Cells.Item(136, 11).Value = 3353.85  # K136: 2189.18922 -> 3353.85
$ws.Cells.Item(136, 12).Value = 8105.571599999999  # L136: 8741.143199999999 -> 8105.571599999999
$ws.Cells.Item(136, 13).Value = -803.8500000000004  # M136: 360.8107799999998 -> -803.8500000000004
$ws.Cells.Item(136, 14).Value = -13205.5716  # N136: -13841.1432 -> -13205.5716

# Sheet BSM, Row 20
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 23260564  # H20: 27032380 -> 23260564
$ws.Cells.Item(20, 9).Value = 38467790  # I20: 58832200 -> 38467790
$ws.Cells.Item(20, 10).Value = 2450.8235  # J20: 2532.75 -> 2450.8235
$ws.Cells.Item(20, 11).Value = 38467790  # K20: 58832200 -> 38467790
$ws.Cells.Item(20, 12).Value = 2450.8235  # L20: 2532.75 -> 2450.8235
$ws.Cells.Item(20, 13).Value = -38467543  # M20: -58831953 -> -38467543
$ws.Cells.Item(20, 14).Value = -2944.8235  # N20: -3026.75 -> -2944.8235

# Sheet BSM, Row 74
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(74, 8).Value = 33880  # H74: 37800 -> 33880
$ws.Cells.Item(74, 10).Value = 33880  # J74: 37800 -> 33880
$ws.Cells.Item(74, 12).Value = 33880  # L74: 37800 -> 33880
$ws.Cells.Item(74, 14).Value = -35752  # N74: -39672 -> -35752

# Sheet BSM, Row 77
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(77, 8).Value = 33880  # H77: 37800 -> 33880
$ws.Cells.Item(77, 10).Value = 33880  # J77: 37800 -> 33880
$ws.Cells.Item(77, 12).Value = 101640  # L77: 113400 -> 101640
$ws.Cells.Item(77, 14).Value = -111000  # N77: -122760 -> -111000

# Sheet BSM, Row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 987.43634  # H134: 1060.0416 -> 987.43634
$ws.Cells.Item(134, 9).Value = 772.5454999999999  # I134: 855.38464 -> 772.5454999999999
$ws.Cells.Item(134, 10).Value = 1847  # J134: 1946.8889 -> 1847
$ws.Cells.Item(134, 11).Value = 2317.6365  # K134: 2566.15392 -> 2317.6365
$ws.Cells.Item(134, 12).Value = 5541  # L134: 5840.6667 -> 5541
$ws.Cells.Item(134, 13).Value = 217.3635000000004  # M134: -31.15391999999974 -> 217.3635000000004
$ws.Cells.Item(134, 14).Value = -10611  # N134: -10910.6667 -> -10611

# Sheet CRP, Row 7
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 105.347824  # H7: 103 -> 105.347824
$ws.Cells.Item(7, 9).Value = 110.46667  # I7: 106.625 -> 110.46667
$ws.Cells.Item(7, 11).Value = 110.46667  # K7: 106.625 -> 110.46667
$ws.Cells.Item(7, 13).Value = 2.533330000000007  # M7: 6.375 -> 2.533330000000007

# Sheet CRP, Row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(132, 8).Value = 1622.8158  # H132: 1818.3334 -> 1622.8158
$ws.Cells.Item(132, 9).Value = 1345.3214  # I132: 1567.2609 -> 1345.3214
$ws.Cells.Item(132, 10).Value = 2399.8  # J132: 2395.8 -> 2399.8
$ws.Cells.Item(132, 11).Value = 4035.9642  # K132: 4701.7827 -> 4035.9642
$ws.Cells.Item(132, 12).Value = 7199.400000000001  # L132: 7187.400000000001 -> 7199.400000000001
$ws.Cells.Item(132, 13).Value = -1505.9642  # M132: -2171.7827 -> -1505.9642
$ws.Cells.Item(132, 14).Value = -12259.4  # N132: -12247.4 -> -12259.4

# Sheet CRP, Row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(134, 8).Value = 1415.7333  # H134: 1379.6154 -> 1415.7333
$ws.Cells.Item(134, 9).Value = 1203.2727  # I134: 1173.5 -> 1203.2727
$ws.Cells.Item(134, 10).Value = 2000  # J134: 2066.6667 -> 2000
$ws.Cells.Item(134, 11).Value = 3609.8181  # K134: 3520.5 -> 3609.8181
$ws.Cells.Item(134, 12).Value = 6000  # L134: 6200.000100000001 -> 6000
$ws.Cells.Item(134, 13).Value = -1074.8181  # M134: -985.5 -> -1074.8181
$ws.Cells.Item(134, 14).Value = -11070  # N134: -11270.0001 -> -11070

# Sheet GSM, Row 70
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 6169.143  # H70: 6634.857 -> 6169.143
$ws.Cells.Item(70, 9).Value = 5688  # I70: 5907.9 -> 5688
$ws.Cells.Item(70, 10).Value = 7933.3335  # J70: 8452.25 -> 7933.3335
$ws.Cells.Item(70, 11).Value = 5688  # K70: 5907.9 -> 5688
$ws.Cells.Item(70, 12).Value = 7933.3335  # L70: 8452.25 -> 7933.3335
$ws.Cells.Item(70, 13).Value = -5418  # M70: -5637.9 -> -5418
$ws.Cells.Item(70, 14).Value = -8473.333500000001  # N70: -8992.25 -> -8473.333500000001

# Sheet GSM, Row 73
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(73, 8).Value = 6169.143  # H73: 6634.857 -> 6169.143
$ws.Cells.Item(73, 9).Value = 5688  # I73: 5907.9 -> 5688
$ws.Cells.Item(73, 10).Value = 7933.3335  # J73: 8452.25 -> 7933.3335
$ws.Cells.Item(73, 11).Value = 5688  # K73: 5907.9 -> 5688
$ws.Cells.Item(73, 12).Value = 7933.3335  # L73: 8452.25 -> 7933.3335
$ws.Cells.Item(73, 13).Value = -4752  # M73: -4971.9 -> -4752
$ws.Cells.Item(73, 14).Value = -9805.333500000001  # N73: -10324.25 -> -9805.333500000001

# Sheet LTW, Row 68
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(68, 8).Value = 8754.134  # H68: 9305 -> 8754.134
$ws.Cells.Item(68, 9).Value = 26403  # I68: 21474 -> 26403
$ws.Cells.Item(68, 10).Value = 2336.3635  # J68: 2544.4443 -> 2336.3635
$ws.Cells.Item(68, 11).Value = 26403  # K68: 21474 -> 26403
$ws.Cells.Item(68, 12).Value = 2336.3635  # L68: 2544.4443 -> 2336.3635
$ws.Cells.Item(68, 13).Value = -25654  # M68: -20725 -> -25654
$ws.Cells.Item(68, 14).Value = -3834.3635  # N68: -4042.4443 -> -3834.3635

# Sheet LTW, Row 71
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(71, 8).Value = 8754.134  # H71: 9305 -> 8754.134
$ws.Cells.Item(71, 9).Value = 26403  # I71: 21474 -> 26403
$ws.Cells.Item(71, 10).Value = 2336.3635  # J71: 2544.4443 -> 2336.3635
$ws.Cells.Item(71, 11).Value = 132015  # K71: 107370 -> 132015
$ws.Cells.Item(71, 12).Value = 11681.8175  # L71: 12722.2215 -> 11681.8175
$ws.Cells.Item(71, 13).Value = -128271  # M71: -103626 -> -128271
$ws.Cells.Item(71, 14).Value = -19169.8175  # N71: -20210.2215 -> -19169.8175

# Sheet LTW, Row 119
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(119, 8).Value = 0  # H119: 27710 -> 0
$ws.Cells.Item(119, 10).Value = 0  # J119: 27710 -> 0
$ws.Cells.Item(119, 12).Value = 0  # L119: 27710 -> 0
$ws.Cells.Item(119, 14).ClearContents()  # N119: was -37386

# Sheet WVR, Row 45
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(45, 8).Value = 6586.857  # H45: 8415.333000000001 -> 6586.857
$ws.Cells.Item(45, 10).Value = 6586.857  # J45: 8415.333000000001 -> 6586.857
$ws.Cells.Item(45, 12).Value = 6586.857  # L45: 8415.333000000001 -> 6586.857
$ws.Cells.Item(45, 14).Value = -7568.857  # N45: -9397.333000000001 -> -7568.857

# Sheet WVR, Row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 1285.3871  # H132: 868.96924 -> 1285.3871
$ws.Cells.Item(132, 9).Value = 1101.9333  # I132: 723.2558 -> 1101.9333
$ws.Cells.Item(132, 10).Value = 1457.375  # J132: 1153.7727 -> 1457.375
$ws.Cells.Item(132, 11).Value = 3305.7999  # K132: 2169.7674 -> 3305.7999
$ws.Cells.Item(132, 12).Value = 4372.125  # L132: 3461.3181 -> 4372.125
$ws.Cells.Item(132, 13).Value = -775.7999  # M132: 360.2325999999998 -> -775.7999
$ws.Cells.Item(132, 14).Value = -9432.125  # N132: -8521.3181 -> -9432.125

# Sheet WVR, Row 133
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(133, 8).Value = 20715  # H133: 30715 -> 20715
$ws.Cells.Item(133, 10).Value = 20715  # J133: 30715 -> 20715
$ws.Cells.Item(133, 12).Value = 20715  # L133: 30715 -> 20715
$ws.Cells.Item(133, 14).Value = -30835  # N133: -40835 -> -30835

# Sheet WVR, Row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(136, 8).Value = 3704.7837  # H136: 3022.6667 -> 3704.7837
$ws.Cells.Item(136, 9).Value = 725.0526  # I136: 525.15625 -> 725.0526
$ws.Cells.Item(136, 10).Value = 6850.0557  # J136: 9170.385 -> 6850.0557
$ws.Cells.Item(136, 11).Value = 2175.1578  # K136: 1575.46875 -> 2175.1578
$ws.Cells.Item(136, 12).Value = 20550.1671  # L136: 27511.155 -> 20550.1671
$ws.Cells.Item(136, 13).Value = 374.8422  # M136: 974.53125 -> 374.8422
$ws.Cells.Item(136, 14).Value = -25650.1671  # N136: -32611.155 -> -25650.1671
